# Updated symbol list on Sat Dec 24 14:41:11 UTC 2022 with GitHub Actions
#
# This script updates price (column D) and volume-label (column E) cells
# on the active worksheet to reflect the latest scraped crypto data.
# Because several of the new values still look like plain numbers (e.g.
# "244.32"), assigning them directly through .Value would make Excel
# auto-convert the cell to a numeric type. The source file stores these
# cells as text, so we briefly force a Text number format before writing
# the value, then restore the cell style to "Normal" (matching the
# original, unstyled cells) so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue "D2"  "244.32"
Set-TextValue "D3"  "21.91"
Set-TextValue "D4"  "5.393"
Set-TextValue "D5"  "0.05985"
Set-TextValue "D6"  "3.391"
Set-TextValue "D7"  "0.8146"
Set-TextValue "D8"  "0.9543"
Set-TextValue "D10" "0.07422"
Set-TextValue "D11" "0.03254"
Set-TextValue "D12" "0.03077"
Set-TextValue "D13" "0.09407"
Set-TextValue "D14" "4.001"
Set-TextValue "D15" "0.001594"
Set-TextValue "D16" "0.04799"
Set-TextValue "D17" "0.0005901"
Set-TextValue "D18" "0.005448"
Set-TextValue "D19" "0.004147"
Set-TextValue "D20" "0.0009885"
Set-TextValue "D22" "3.680"
Set-TextValue "D23" "6.428"
Set-TextValue "D24" "2.188"
Set-TextValue "D40" "0.04003"
Set-TextValue "D41" "0.006651"
Set-TextValue "D42" "0.1073"
Set-TextValue "D44" "0.005732"
Set-TextValue "D45" "0.00005124"
Set-TextValue "D48" "0.006658"
Set-TextValue "D49" "0.00002100"

# Column E (Volume(1h) label) updates
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
